$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values look like plain numbers (e.g. "301.99").
# Excel would otherwise auto-convert them to numeric values, so we
# temporarily force a Text number format, assign the value, then
# restore the default "Normal" style so the cell keeps its original look.
$textCells = @('D5', 'D6', 'D7', 'D10', 'D11', 'D12', 'D19', 'D22', 'D23', 'D27', 'D28', 'D31', 'D33', 'D35', 'D36', 'D39', 'D44', 'D45', 'D46', 'D48')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D5').Value = '301.99'
$ws.Range('D6').Value = '95.88'
$ws.Range('D7').Value = '0.504'
$ws.Range('D10').Value = '34.72'
$ws.Range('D11').Value = '0.0780'
$ws.Range('D12').Value = '18.54'
$ws.Range('D19').Value = '12.80'
$ws.Range('D22').Value = '67.03'
$ws.Range('D23').Value = '236.01'
$ws.Range('D27').Value = '24.57'
$ws.Range('D28').Value = '167.56'
$ws.Range('D31').Value = '32.77'
$ws.Range('D33').Value = '17.78'
$ws.Range('D35').Value = '4.45'
$ws.Range('D36').Value = '2.36'
$ws.Range('D39').Value = '1.73'
$ws.Range('D44').Value = '10.15'
$ws.Range('D45').Value = '18.25'
$ws.Range('D46').Value = '2.15'
$ws.Range('D48').Value = '53.44'

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}

# Remaining updates (plain text, not at risk of numeric auto-conversion)
$ws.Range('D2').Value = '42.654.07'
$ws.Range('E2').Value = '  -0.64%  '
$ws.Range('D3').Value = '2.294.21'
$ws.Range('E3').Value = '  -0.37%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('E5').Value = '  +0.69%  '
$ws.Range('E6').Value = '  -1.82%  '
$ws.Range('E7').Value = '  -2.02%  '
$ws.Range('E9').Value = '  -2.24%  '
$ws.Range('E10').Value = '  -3.13%  '
$ws.Range('E11').Value = '  -0.87%  '
$ws.Range('E12').Value = '  +4.70%  '
$ws.Range('E13').Value = '  +2.29%  '
$ws.Range('E14').Value = '  +0.93%  '
$ws.Range('D15').Value = '2.651.82'
$ws.Range('E15').Value = '  -0.30%  '
$ws.Range('D16').Value = '2.279.38'
$ws.Range('E16').Value = '  -0.94%  '
$ws.Range('E17').Value = '  -0.86%  '
$ws.Range('D18').Value = '42.585.49'
$ws.Range('E18').Value = '  -0.72%  '
$ws.Range('E19').Value = '  +1.59%  '
$ws.Range('D20').Value = '0.0₃0890'
$ws.Range('E20').Value = '  -1.83%  '
$ws.Range('E21').Value = '  -1.81%  '
$ws.Range('E22').Value = '  -1.40%  '
$ws.Range('E23').Value = '  -2.60%  '
$ws.Range('E24').Value = '  -0.56%  '
$ws.Range('E25').Value = '  +0.20%  '
$ws.Range('E26').Value = '  -1.62%  '
$ws.Range('E27').Value = '  -1.93%  '
$ws.Range('E28').Value = '  +0.74%  '
$ws.Range('E30').Value = '  -0.60%  '
$ws.Range('E31').Value = '  +0.10%  '
$ws.Range('E32').Value = '  +0.07%  '
$ws.Range('E33').Value = '  +1.69%  '
$ws.Range('E34').Value = '  -1.15%  '
$ws.Range('E35').Value = '  -6.34%  '
$ws.Range('E36').Value = '  -1.78%  '
$ws.Range('E37').Value = '  -0.17%  '
$ws.Range('E38').Value = '  -0.55%  '
$ws.Range('E39').Value = '  -1.77%  '
$ws.Range('E40').Value = '  -1.11%  '
$ws.Range('E41').Value = '  -2.68%  '
$ws.Range('D42').Value = '1.991.81'
$ws.Range('E42').Value = '  -0.38%  '
$ws.Range('E43').Value = '  -2.16%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('E44').Value = '  +0.15%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('E45').Value = '  +5.62%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('E46').Value = '  -0.19%  '
$ws.Range('E47').Value = '  -0.72%  '
$ws.Range('E48').Value = '  +0.17%  '
$ws.Range('E49').Value = '  +4.55%  '
$ws.Range('D50').Value = '2.518.71'
$ws.Range('E51').Value = '  -2.14%  '
